# "added cards to members app"
# - Renumber the image filenames for rows 11-21 (member cards 10-20) from
#   "upload/00NN.JPG" (extra leading zero) to "upload/0NN.JPG".
# - Normalize the "Midfielder" area-of-expertise label to "Midfield" so it
#   matches the rest of the lookup table (rows 2, 7, 12, 17).
# - Re-apply the sheet's remembered filter range (hidden _FilterDatabase name).
# - Leave the sheet's selection where the author left off (cell B9).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Memberlist")

# Fix the image card filenames (drop the extra leading zero) for rows 11-21.
$imageFixes = @{
    11 = "upload/010.JPG"
    12 = "upload/011.JPG"
    13 = "upload/012.JPG"
    14 = "upload/013.JPG"
    15 = "upload/014.JPG"
    16 = "upload/015.JPG"
    17 = "upload/016.JPG"
    18 = "upload/017.JPG"
    19 = "upload/018.JPG"
    20 = "upload/019.JPG"
    21 = "upload/020.JPG"
}
foreach ($row in $imageFixes.Keys) {
    $ws.Range("B$row").Value = $imageFixes[$row]
}

# Normalize "Midfielder" -> "Midfield" in the area_of_expertise column.
$midfieldRows = @(2, 7, 12, 17)
foreach ($row in $midfieldRows) {
    $ws.Range("H$row").Value = "Midfield"
}

# Restore the hidden AutoFilter range defined name (no active filter arrows).
$fdb = $ws.Names.Add("_xlnm._FilterDatabase", "=Memberlist!`$A`$1:`$N`$21")
$fdb.Visible = $false

# Move the active selection to B9, as left by the author.
[void]$ws.Range("B9").Select()
